# Trade #18 closed at 2026-02-16 22:59:08 - base_strategy DOWN +0.000%
#
# The "All Trades" and "base_strategy" sheets both keep a mirrored trade
# log; append the newly-opened trade (row 19) to each, copying the
# existing layout of the last row (row 18).

$wb  = $excel.ActiveWorkbook
$row = 19
$prevRow = 18

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Exit Price (G) and Exit Reason (P) stay blank for an OPEN trade.
    # Copy the blank template cells from the previous row so the new row
    # keeps the same blank-cell layout as every other row in the sheet.
    $ws.Range("G" + $prevRow).Copy($ws.Range("G" + $row))
    $ws.Range("P" + $prevRow).Copy($ws.Range("P" + $row))

    $ws.Cells.Item($row, 1).Value = 18

    # Force the date to stay literal text instead of Excel's automatic
    # date recognition (this workbook stores dates as plain text, like
    # every other row already on the sheet).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value  = "22:59:08"
    $ws.Cells.Item($row, 4).Value  = "base_strategy"
    $ws.Cells.Item($row, 5).Value  = "DOWN"
    $ws.Cells.Item($row, 6).Value  = 0.5
    $ws.Cells.Item($row, 8).Value  = "OPEN"
    $ws.Cells.Item($row, 9).Value  = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 17).Value = 0
}
